# Generate Report for Handback
# Update handback timestamps for the "5a860db5-99f2-4943-9fd0-e75c3d37a981" entry
# across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 00:45:40"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-08-25 00:45:36"
$wsZh.Range("K2").Value = "2016-08-25 00:45:53"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-08-25 00:45:40"
$wsDe.Range("K2").Value = "2016-08-25 00:46:01"
